# Insert a new data row at spreadsheet row 697 (pushing existing rows
# 697..788 down to 698..789) and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(697).Insert()

$ws.Range("A697").Value = 6
$ws.Range("B697").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C697").Value = "Metropolitana"
$ws.Range("D697").Value = 45127
$ws.Range("E697").Value = 13
$ws.Range("F697").Value = 100112012
$ws.Range("G697").Value = "Espinaca"
$ws.Range("H697").Value = "Sin especificar"
$ws.Range("I697").Value = "Primera"
$ws.Range("J697").Value = 550
$ws.Range("K697").Value = 5000
$ws.Range("L697").Value = 6000
$ws.Range("M697").Value = 5455
$ws.Range("N697").Value = "$/cuna 10 kilos"
$ws.Range("O697").Value = "Región Metropolitana"
$ws.Range("P697").Value = 546
$ws.Range("Q697").Value = 10
$ws.Range("R697").Value = "Hortaliza"
